# Insert a new weekly price record as row 51 (pushing the existing
# rows 51-81 down to 52-82), matching the "Fruta / hortaliza, semanal"
# update for Berenjena @ Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 51..81 down by one (EntireRow insert keeps their content/style).
$ws.Rows.Item(51).Insert()

# Populate the newly-inserted row 51 with the new weekly observation.
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 45072
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 100112001
$ws.Range("G51").Value = "Berenjena"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 80
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 7000
$ws.Range("M51").Value = 6500
$ws.Range("N51").Value = "`$/caja 60 unidades"
$ws.Range("O51").Value = "Región de Arica y Parinacota"
$ws.Range("P51").Value = 108
$ws.Range("Q51").Value = 60
$ws.Range("R51").Value = "Hortaliza"
